$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value in A2 (creates a new shared string "update")
$ws.Range("A2").Value = "update"

# Make A2 the active/selected cell, matching the saved selection in the sheet view
$ws.Range("A2").Select()
